# Update TestData.xlsx with latest ParaBank test data
#
# The "ParaBank_RegistrationForm" sheet's row-2 result message (column M)
# is rewritten to wrap across three lines, its status (column N) flips
# from FAIL to PASS, and the message cell gets wrap-text formatting so the
# new line breaks render correctly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ParaBank_RegistrationForm")

# Column M row 2: multi-line success message (replaces the old single-line text)
$ws.Range("M2").Value = "Your account was created" + [char]10 + "successfully. You are now" + [char]10 + "logged in."
$ws.Range("M2").WrapText = $true

# Column N row 2: status flips from FAIL to PASS
$ws.Range("N2").Value = "PASS"
